$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.749.65"
$ws.Range("E2").Value = "  -3.34%  "

$ws.Range("D3").Value = "2.983.09"
$ws.Range("E3").Value = "  -2.62%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.61"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.03"
$ws.Range("E6").Value = "  -2.03%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "2.977.71"
$ws.Range("E8").Value = "  -2.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -1.44%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("E10").Value = "  -4.08%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.12"
$ws.Range("E11").Value = "  -1.79%  "

$ws.Range("E12").Value = "  -3.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000218"
$ws.Range("E13").Value = "  -2.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.44"
$ws.Range("E14").Value = "  -2.73%  "

$ws.Range("D15").Value = "3.464.01"
$ws.Range("E15").Value = "  -2.57%  "

$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "60.858.51"
$ws.Range("E17").Value = "  -3.11%  "

$ws.Range("D18").Value = "2.984.27"
$ws.Range("E18").Value = "  -2.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "458.28"
$ws.Range("E20").Value = "  -5.39%  "

$ws.Range("E21").Value = "  -1.48%  "

$ws.Range("E22").Value = "  -3.92%  "

$ws.Range("E23").Value = "  -4.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.40"
$ws.Range("E24").Value = "  -1.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.78"
$ws.Range("E25").Value = "  -3.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  -2.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.69"
$ws.Range("E28").Value = "  -5.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.14"
$ws.Range("E30").Value = "  +2.25%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.85"
$ws.Range("E31").Value = "  -1.28%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.18"
$ws.Range("E32").Value = "  -3.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "55.03"
$ws.Range("E33").Value = "  -2.49%  "

$ws.Range("E34").Value = "  -5.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("E35").Value = "  -1.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.77"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "454.92"
$ws.Range("E37").Value = "  -4.75%  "

$ws.Range("D38").Value = "3.160.41"
$ws.Range("E38").Value = "  +2.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0778"
$ws.Range("E39").Value = "  -2.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0379"
$ws.Range("E40").Value = "  -4.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.116"
$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.04"
$ws.Range("E42").Value = "  -0.66%  "

$ws.Range("E43").Value = "  -9.71%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.93"
$ws.Range("E45").Value = "  +5.03%  "

$ws.Range("E46").Value = "  -4.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.57"
$ws.Range("E47").Value = "  -1.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.107"
$ws.Range("E48").Value = "  -0.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.96"
$ws.Range("E49").Value = "  -4.12%  "

$ws.Range("D50").Value = "0.0₃0491"
$ws.Range("E50").Value = "  -10.67%  "

$ws.Range("E51").Value = "  +7.72%  "
